# Fix product importing: add missing "line total" formulas in column L
# (quantity * unit price) for the two data rows, and update the saved
# view state (zoom level + active selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing formulas to L4 and L5 (J*K = quantity * price)
$ws.Range("L4").Formula = "=J4*K4"
$ws.Range("L5").Formula = "=J5*K5"

# Update the view: zoom to 140% and move the selection to L5
$excel.ActiveWindow.Zoom = 140
$ws.Range("L5").Select()
